$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 158.4
$ws.Range("I11").Value = 158.4
$ws.Range("K11").Value = 158.4
$ws.Range("M11").Value = -18.40000000000001

$ws.Range("H33").Value = 573.8125
$ws.Range("I33").Value = 569.4545000000001
$ws.Range("K33").Value = 569.4545000000001
$ws.Range("M33").Value = -340.4545000000001

$ws.Range("H34").Value = 1499.3334
$ws.Range("J34").Value = 1499
$ws.Range("L34").Value = 1499
$ws.Range("N34").Value = -1905

$ws.Range("H36").Value = 1499.3334
$ws.Range("J36").Value = 1499
$ws.Range("L36").Value = 1499
$ws.Range("N36").Value = -2929

$ws.Range("H40").Value = 4149.9
$ws.Range("I40").Value = 3166.6667
$ws.Range("J40").Value = 5624.75
$ws.Range("K40").Value = 3166.6667
$ws.Range("L40").Value = 5624.75
$ws.Range("M40").Value = -2991.6667
$ws.Range("N40").Value = -5974.75

$ws.Range("H51").Value = 10399
$ws.Range("J51").Value = 10399
$ws.Range("L51").Value = 10399
$ws.Range("N51").Value = -11367

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H86").Value = 3572
$ws.Range("I86").Value = 3572
$ws.Range("K86").Value = 3572
$ws.Range("M86").Value = -2449

$ws.Range("H89").Value = 3572
$ws.Range("I89").Value = 3572
$ws.Range("K89").Value = 17860
$ws.Range("M89").Value = -12244

$ws.Range("H98").Value = 628.2222
$ws.Range("I98").Value = 628.2222
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 628.2222
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 869.7778
$ws.Range("N98").ClearContents()

$ws.Range("H106").Value = 18580.6
$ws.Range("I106").Value = 18580.6
$ws.Range("K106").Value = 18580.6
$ws.Range("M106").Value = -17949.6

$ws.Range("H107").Value = 967.5714
$ws.Range("I107").Value = 678.8333
$ws.Range("K107").Value = 678.8333
$ws.Range("M107").Value = 1241.1667

$ws.Range("H122").Value = 628.2222
$ws.Range("I122").Value = 628.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1884.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 565.3334
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 1248
$ws.Range("I125").Value = 1100
$ws.Range("K125").Value = 9900
$ws.Range("M125").Value = -7440

$ws.Range("H131").Value = 6841.857
$ws.Range("I131").Value = 2633.3333
$ws.Range("J131").Value = 9998.25
$ws.Range("K131").Value = 7899.999899999999
$ws.Range("L131").Value = 29994.75
$ws.Range("M131").Value = -2859.999899999999
$ws.Range("N131").Value = -40074.75

$ws.Range("H132").Value = 32261436
$ws.Range("I132").Value = 32261436
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 96784308
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -96781778
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 1211.5454
$ws.Range("I135").Value = 1211.5454
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10903.9086
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -8368.908599999999
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1395.8334
$ws.Range("I45").Value = 1395.8334
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1395.8334
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1018.8334
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 1909.375
$ws.Range("I61").Value = 1897.1428
$ws.Range("K61").Value = 1897.1428
$ws.Range("M61").Value = -1685.1428

$ws.Range("H102").Value = 2633.3333
$ws.Range("I102").Value = 2633.3333
$ws.Range("K102").Value = 2633.3333
$ws.Range("M102").Value = -1011.3333

$ws.Range("H136").Value = 1909.375
$ws.Range("I136").Value = 1897.1428
$ws.Range("K136").Value = 5691.428400000001
$ws.Range("M136").Value = -3141.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 11000
$ws.Range("J8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("N8").Value = -11280

$ws.Range("H99").Value = 3666.6667
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2

$ws.Range("H105").Value = 4048.5
$ws.Range("I105").Value = 4198.2856
$ws.Range("K105").Value = 4198.2856
$ws.Range("M105").Value = -2451.2856

$ws.Range("H134").Value = 1896.9
$ws.Range("J134").Value = 2000
$ws.Range("L134").Value = 6000
$ws.Range("N134").Value = -11070

$ws.Range("H141").Value = 76943
$ws.Range("J141").Value = 76925.664
$ws.Range("L141").Value = 76925.664
$ws.Range("N141").Value = -87285.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7294.7827
$ws.Range("J31").Value = 5246.643
$ws.Range("L31").Value = 5246.643
$ws.Range("N31").Value = -5836.643

$ws.Range("H34").Value = 7294.7827
$ws.Range("J34").Value = 5246.643
$ws.Range("L34").Value = 5246.643
$ws.Range("N34").Value = -5650.643

$ws.Range("H132").Value = 10008513
$ws.Range("I132").Value = 15395408
$ws.Range("K132").Value = 46186224
$ws.Range("M132").Value = -46183694

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 113888.78
$ws.Range("I5").Value = 3124.25
$ws.Range("K5").Value = 9372.75
$ws.Range("M5").Value = -9260.75

$ws.Range("H12").Value = 34.583332
$ws.Range("I12").Value = 9.25
$ws.Range("J12").Value = 47.25
$ws.Range("K12").Value = 27.75
$ws.Range("L12").Value = 141.75
$ws.Range("M12").Value = 145.25
$ws.Range("N12").Value = -487.75

$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H135").Value = 113888.78
$ws.Range("I135").Value = 3124.25
$ws.Range("K135").Value = 28118.25
$ws.Range("M135").Value = -25583.25

$ws.Range("H138").Value = 9812.412
$ws.Range("I138").Value = 8214.643
$ws.Range("K138").Value = 24643.929
$ws.Range("M138").Value = -19503.929

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 1457.1428
$ws.Range("J59").Value = 3000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4166

$ws.Range("H70").Value = 3750
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 3750
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H97").Value = 1042.7142
$ws.Range("I97").Value = 733.1667
$ws.Range("K97").Value = 733.1667
$ws.Range("M97").Value = -237.1667

$ws.Range("H107").Value = 3231.5557
$ws.Range("J107").Value = 6499.25
$ws.Range("L107").Value = 6499.25
$ws.Range("N107").Value = -10339.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1106.1428
$ws.Range("I16").Value = 978.1667
$ws.Range("K16").Value = 978.1667
$ws.Range("M16").Value = -808.1667

$ws.Range("H40").Value = 31253478
$ws.Range("I40").Value = 45457460
$ws.Range("J40").Value = 4720.8
$ws.Range("K40").Value = 45457460
$ws.Range("L40").Value = 4720.8
$ws.Range("M40").Value = -45457324
$ws.Range("N40").Value = -4992.8

$ws.Range("H46").Value = 3711.1177
$ws.Range("I46").Value = 1760
$ws.Range("J46").Value = 4775.364
$ws.Range("K46").Value = 1760
$ws.Range("L46").Value = 4775.364
$ws.Range("M46").Value = -1572
$ws.Range("N46").Value = -5151.364

$ws.Range("H98").Value = 59756.4
$ws.Range("J98").Value = 59756.4
$ws.Range("L98").Value = 59756.4
$ws.Range("N98").Value = -65746.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 14907.8
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 14907.8
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 14907.8
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16405.8

$ws.Range("H72").Value = 14907.8
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 14907.8
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 44723.39999999999
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -52211.39999999999

$ws.Range("H113").Value = 245.76471
$ws.Range("I113").Value = 212.26666
$ws.Range("K113").Value = 636.79998
$ws.Range("M113").Value = 1533.20002
